$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '27.393.44'
$ws.Cells.Item(2, 5).Value = '  -0.72%  '
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.818.16'
$ws.Cells.Item(3, 5).Value = '  -2.25%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.003'
$ws.Cells.Item(4, 5).Value = '  -1.18%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '332.33'
$ws.Cells.Item(5, 5).Value = '  -0.82%  '
$ws.Cells.Item(6, 5).Value = '  -1.13%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4555'
$ws.Cells.Item(7, 5).Value = '  -1.98%  '
$ws.Cells.Item(8, 5).Value = '  -2.70%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '45.86'
$ws.Cells.Item(9, 5).Value = '  -0.27%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.07822'
$ws.Cells.Item(10, 5).Value = '  -1.62%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.9543'
$ws.Cells.Item(11, 5).Value = '  -4.34%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '20.90'
$ws.Cells.Item(12, 5).Value = '  -3.33%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '1.812.98'
$ws.Cells.Item(13, 5).Value = '  -3.08%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '5.800'
$ws.Cells.Item(14, 5).Value = '  -2.17%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '7.021'
$ws.Cells.Item(15, 5).Value = '  -2.34%  '
$ws.Cells.Item(16, 5).Value = '  -1.03%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '89.09'
$ws.Cells.Item(17, 5).Value = '  +1.04%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.06556'
$ws.Cells.Item(18, 5).Value = '  -2.58%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.00001015'
$ws.Cells.Item(19, 5).Value = '  -2.56%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '17.00'
$ws.Cells.Item(20, 5).Value = '  -1.03%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '1.002'
$ws.Cells.Item(21, 5).Value = '  -1.17%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '27.381.06'
$ws.Cells.Item(22, 5).Value = '  -0.79%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '5.257'
$ws.Cells.Item(23, 5).Value = '  -3.51%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '10.76'
$ws.Cells.Item(24, 5).Value = '  -1.54%  '
$ws.Cells.Item(25, 5).Value = '  -2.19%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '158.10'
$ws.Cells.Item(26, 5).Value = '  -0.96%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '2.024.42'
$ws.Cells.Item(27, 5).Value = '  -3.10%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '19.21'
$ws.Cells.Item(28, 5).Value = '  -2.04%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.029'
$ws.Cells.Item(29, 5).Value = '  -5.23%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '5.242'
$ws.Cells.Item(30, 5).Value = '  -3.13%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '117.42'
$ws.Cells.Item(31, 5).Value = '  -3.46%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.09305'
$ws.Cells.Item(32, 5).Value = '  -1.33%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.9263'
$ws.Cells.Item(33, 5).Value = '  -4.86%  '
$ws.Cells.Item(34, 5).Value = '  -1.80%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '5.194'
$ws.Cells.Item(35, 5).Value = '  -2.04%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.306'
$ws.Cells.Item(36, 5).Value = '  -2.32%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.05880'
$ws.Cells.Item(37, 5).Value = '  -2.26%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.02178'
$ws.Cells.Item(38, 5).Value = '  -2.35%  '
$ws.Cells.Item(39, 2).Value = 'Frax'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.002'
$ws.Cells.Item(39, 5).Value = '  -1.15%  '
$ws.Cells.Item(40, 2).Value = 'FraxShare'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '8.044'
$ws.Cells.Item(40, 5).Value = '  -3.26%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '1.137'
$ws.Cells.Item(41, 5).Value = '  -4.71%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.5706'
$ws.Cells.Item(42, 5).Value = '  -3.87%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.1805'
$ws.Cells.Item(43, 5).Value = '  -3.24%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '9.853'
$ws.Cells.Item(44, 5).Value = '  -4.27%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.261'
$ws.Cells.Item(45, 5).Value = '  +0.72%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.5366'
$ws.Cells.Item(46, 5).Value = '  -3.89%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '11.70'
$ws.Cells.Item(47, 5).Value = '  -3.47%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.862'
$ws.Cells.Item(48, 5).Value = '  -2.91%  '
$ws.Cells.Item(49, 2).Value = 'Cronos'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.06559'
$ws.Cells.Item(49, 5).Value = '  -2.61%  '
$ws.Cells.Item(50, 2).Value = 'Quant'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '109.79'
$ws.Cells.Item(50, 5).Value = '  -1.38%  '
$ws.Cells.Item(51, 5).Value = '  -33.19%  '
